$d = $word.ActiveDocument

$replacements = @(
    @{Old = "171×9=1539"; New = "771×4=3084"},
    @{Old = "982×9=8838"; New = "592×6=3552"},
    @{Old = "754×2=1508"; New = "525×4=2100"},
    @{Old = "949×5=4745"; New = "865×2=1730"},
    @{Old = "821×2=1642"; New = "835×7=5845"},
    @{Old = "375×5=1875"; New = "306×2=612"},
    @{Old = "799×7=5593"; New = "372×9=3348"},
    @{Old = "628×4=2512"; New = "976×4=3904"},
    @{Old = "980×8=7840"; New = "736×2=1472"},
    @{Old = "323×8=2584"; New = "865×4=3460"},
    @{Old = "543×2=1086"; New = "221×3=663"},
    @{Old = "949×3=2847"; New = "878×9=7902"},
    @{Old = "422×3=1266"; New = "354×8=2832"},
    @{Old = "851×7=5957"; New = "500×5=2500"},
    @{Old = "406×3=1218"; New = "294×4=1176"},
    @{Old = "894×8=7152"; New = "395×6=2370"},
    @{Old = "405×9=3645"; New = "127×3=381"},
    @{Old = "736×9=6624"; New = "300×5=1500"},
    @{Old = "191×5=955"; New = "199×5=995"},
    @{Old = "319×8=2552"; New = "407×6=2442"},
    @{Old = "778×6=4668"; New = "233×5=1165"},
    @{Old = "315×6=1890"; New = "382×4=1528"},
    @{Old = "898×3=2694"; New = "685×3=2055"},
    @{Old = "711×6=4266"; New = "176×5=880"},
    @{Old = "768×6=4608"; New = "404×7=2828"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.New, 2) | Out-Null
}
